# Apply scheduled-runner price/profit updates to the Leve profit tracker.
# Each sheet tracks per-leve market data (H-N) for Final Fantasy XIV crafting leves.
# Values below are literal (no formulas in this workbook) so we just overwrite them.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 7: The Bleat Is On | Maple Wand
$ws.Range("H7").Value = 8000
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

# Row 14: Wand-full Tonight | Budding Maple Wand
$ws.Range("H14").Value = 8000
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

# Row 19: Unbreak My Heart | Roof Tile
$ws.Range("H19").Value = 464.4
$ws.Range("J19").Value = 539.25
$ws.Range("L19").Value = 539.25
$ws.Range("N19").Value = -889.25

# Row 37: The Wailers' First Law of Potion | Hi-Potion
$ws.Range("H37").Value = 1000
$ws.Range("J37").Value = 1000
$ws.Range("L37").Value = 3000
$ws.Range("N37").Value = -3252

# Row 50: A Patch-up Place | Mega-Potion
$ws.Range("H50").Value = 999.5
$ws.Range("J50").Value = 999.5
$ws.Range("L50").Value = 2998.5
$ws.Range("N50").Value = -3948.5

# Row 64: Forged from the Void | Void Glue
$ws.Range("H64").Value = 6400
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 6400
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 6400
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -6896

# Row 67: Dodging the Draft (L) | Void Glue
$ws.Range("H67").Value = 6400
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 6400
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 6400
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -8116


$ws = $wb.Worksheets.Item("ARM")
# Row 4: Eyes Bigger than the Plate | Bronze Plate
$ws.Range("H4").Value = 206
$ws.Range("I4").Value = 48.2
$ws.Range("K4").Value = 48.2
$ws.Range("M4").Value = 67.8

# Row 15: All Ovo That | Iron Skillet
$ws.Range("H15").Value = 19999.5
$ws.Range("J15").Value = 19999.5
$ws.Range("L15").Value = 19999.5
$ws.Range("N15").Value = -20699.5

# Row 17: Cook Intentions | Amateur's Skillet
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

# Row 36: Hot for Teacher | Heavy Iron Armor
$ws.Range("H36").Value = 3180.5
$ws.Range("I36").Value = 2907.3333
$ws.Range("K36").Value = 2907.3333
$ws.Range("M36").Value = -2561.3333


$ws = $wb.Worksheets.Item("BSM")
# Row 7: Thank You for Your Business | Bronze Bastard Sword
$ws.Range("H7").Value = 227
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()


$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent | Maple Lumber
$ws.Range("H7").Value = 244.4762
$ws.Range("I7").Value = 274.26666
$ws.Range("K7").Value = 274.26666
$ws.Range("M7").Value = -161.26666

# Row 22: Driving Up the Wall | Elm Lumber
$ws.Range("H22").Value = 859
$ws.Range("I22").Value = 948.75
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 948.75
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -598.75
$ws.Range("N22").Value = -1200

# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 6279.923
$ws.Range("I31").Value = 2044.4286
$ws.Range("K31").Value = 2044.4286
$ws.Range("M31").Value = -1749.4286

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 6279.923
$ws.Range("I34").Value = 2044.4286
$ws.Range("K34").Value = 2044.4286
$ws.Range("M34").Value = -1842.4286


$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap | Maple Syrup
$ws.Range("H5").Value = 1023.2222
$ws.Range("J5").Value = 1294.1428
$ws.Range("L5").Value = 3882.4284
$ws.Range("N5").Value = -4106.428400000001

# Row 9: Jack of All Plates | Jack-o'-lantern
$ws.Range("H9").Value = 842.5
$ws.Range("J9").Value = 842.5
$ws.Range("L9").Value = 2527.5
$ws.Range("N9").Value = -2975.5

# Row 19: The Bango Zango Diet | Parsnip Salad
$ws.Range("H19").Value = 4332
$ws.Range("J19").Value = 3998.5
$ws.Range("L19").Value = 11995.5
$ws.Range("N19").Value = -12343.5

# Row 24: Rustic Repast | Chicken and Mushrooms
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

# Row 29: For Crumbs' Sake | Honey Muffin
$ws.Range("H29").Value = 45.5

# Row 36: Love's Crumpets Lost | Crumpet
$ws.Range("H36").Value = 1876.5
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 1876.5
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 5629.5
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -5967.5

# Row 43: Sole Survivor | Baked Sole
$ws.Range("H43").Value = 2
$ws.Range("I43").Value = 2
$ws.Range("K43").Value = 6
$ws.Range("M43").Value = 108

# Row 50: Moving Up in the World | Rolanberry Cheese
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()

# Row 53: Rolanberry Fields Forever | Rolanberry Cheese
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()

# Row 60: Drinking to Your Health | Mulled Tea
$ws.Range("H60").Value = 315.5
$ws.Range("I60").Value = 305
$ws.Range("J60").Value = 326
$ws.Range("K60").Value = 915
$ws.Range("L60").Value = 978
$ws.Range("M60").Value = -664
$ws.Range("N60").Value = -1480

# Row 61: Red Letter Day | Rolanberry Lassi
$ws.Range("H61").Value = 200
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()

# Row 92: Oh No Udon | Gyr Abanian Flour
$ws.Range("H92").Value = 298.75
$ws.Range("I92").Value = 298.75
$ws.Range("K92").Value = 896.25
$ws.Range("M92").Value = 351.75

# Row 107: Slippery Service | Frantoio Oil
$ws.Range("H107").Value = 2816.5
$ws.Range("I107").Value = 2379.8
$ws.Range("K107").Value = 7139.400000000001
$ws.Range("M107").Value = -5219.400000000001

# Row 135: Not-so-secret Ingredient | Royal Maple Syrup
$ws.Range("H135").Value = 1023.2222
$ws.Range("J135").Value = 1294.1428
$ws.Range("L135").Value = 11647.2852
$ws.Range("N135").Value = -16717.2852


$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers | Copper Ingot
$ws.Range("H2").Value = 65.125
$ws.Range("I2").Value = 70.8
$ws.Range("J2").Value = 55.666668
$ws.Range("K2").Value = 70.8
$ws.Range("L2").Value = 55.666668
$ws.Range("M2").Value = 42.2
$ws.Range("N2").Value = -281.666668

# Row 113: Copious Crystal Cannons | Manasilver Nugget
$ws.Range("H113").Value = 4428.2856
$ws.Range("I113").Value = 4874.75
$ws.Range("J113").Value = 3833
$ws.Range("K113").Value = 4874.75
$ws.Range("L113").Value = 3833
$ws.Range("M113").Value = -2704.75
$ws.Range("N113").Value = -8173


$ws = $wb.Worksheets.Item("LTW")
# Row 2: Red in the Head | Leather Calot
$ws.Range("H2").Value = 35833.082
$ws.Range("J2").Value = 71999.39999999999
$ws.Range("L2").Value = 71999.39999999999
$ws.Range("N2").Value = -72223.39999999999

# Row 9: From the Sands to the Stage | Leather Himantes
$ws.Range("H9").Value = 467.5
$ws.Range("I9").Value = 371.6
$ws.Range("J9").Value = 947
$ws.Range("K9").Value = 371.6
$ws.Range("L9").Value = 947
$ws.Range("M9").Value = -147.6
$ws.Range("N9").Value = -1395

# Row 12: A Place to Call Helm | Hard Leather Pot Helm
$ws.Range("H12").Value = 1494.5
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

# Row 22: Skin off Their Backs | Aldgoat Leather
$ws.Range("H22").Value = 711.875
$ws.Range("I22").Value = 431.66666
$ws.Range("J22").Value = 880
$ws.Range("K22").Value = 431.66666
$ws.Range("L22").Value = 880
$ws.Range("M22").Value = -136.66666
$ws.Range("N22").Value = -1470

# Row 27: Fire and Hide | Aldgoat Leather
$ws.Range("H27").Value = 711.875
$ws.Range("I27").Value = 431.66666
$ws.Range("J27").Value = 880
$ws.Range("K27").Value = 431.66666
$ws.Range("L27").Value = 880
$ws.Range("M27").Value = -324.66666
$ws.Range("N27").Value = -1094

# Row 31: Open to Attack | Goatskin Jacket
$ws.Range("H31").Value = 257.5
$ws.Range("I31").Value = 15
$ws.Range("J31").Value = 500
$ws.Range("K31").Value = 15
$ws.Range("L31").Value = 500
$ws.Range("M31").Value = 233
$ws.Range("N31").Value = -996

# Row 58: Handle with Care | Peisteskin Cesti
$ws.Range("H58").Value = 1825
$ws.Range("J58").Value = 3500
$ws.Range("L58").Value = 3500
$ws.Range("N58").Value = -4020


$ws = $wb.Worksheets.Item("WVR")
# Row 23: Pants Are Not Enough | Padded Hempen Trousers
$ws.Range("H23").Value = 266.33334
$ws.Range("I23").Value = 266.33334
$ws.Range("K23").Value = 266.33334
$ws.Range("M23").Value = -37.33334000000002

# Row 55: A Matter of Import | Woolen Hat
$ws.Range("H55").Value = 5026
$ws.Range("I55").Value = 2539
$ws.Range("K55").Value = 2539
$ws.Range("M55").Value = -2262

# Row 61: Bundle Up, It's Odd out There | Woolen Deerstalker
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()

